function Replace-Next {
    param($d, $searchText, $replaceText, $pos)
    $r = $d.Range($pos, $d.Content.End)
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return $pos
    }
    $r.Text = $replaceText
    return $r.End
}

$d = $word.ActiveDocument
$pos = 0
$pos = Replace-Next $d "(Verse 1)" "Verse 1:" $pos
$pos = Replace-Next $d "There's a guy named Sajeev" "Yaseen, my Pakistani friend" $pos
$pos = Replace-Next $d "He's always down for a laugh" "With his warm smile that never ends" $pos
$pos = Replace-Next $d "With a love for Nasi Goreng" "From Lahore to Karachi, he's a true gem" $pos
$pos = Replace-Next $d "He's got quite the culinary craft" "Always there for me till the very end" $pos
$pos = Replace-Next $d "(Chorus)" "Chorus:" $pos
$pos = Replace-Next $d "Sajeev, oh Sajeev" "Yaseen, Yaseen" $pos
$pos = Replace-Next $d "He's the life of the party" "A friend so true and keen" $pos
$pos = Replace-Next $d "With a plate of Nasi Goreng" "In his heart, love is seen" $pos
$pos = Replace-Next $d "He'll make you feel hearty" "Yaseen, Yaseen" $pos
$pos = Replace-Next $d "(Verse 2)" "Verse 2:" $pos
$pos = Replace-Next $d "From Kuala Lumpur to Penang" "I remember the days we spent" $pos
$pos = Replace-Next $d "His love for food will never end" "Laughing and talking 'til the night was spent" $pos
$pos = Replace-Next $d "He's always up for a food adventure" "His kindness and wisdom, a guiding light" $pos
$pos = Replace-Next $d "His appetite, he'll always defend" "In his presence, everything feels right" $pos
$pos = Replace-Next $d "(Chorus)" "Chorus:" $pos
$pos = Replace-Next $d "Sajeev, oh Sajeev" "Yaseen, Yaseen" $pos
$pos = Replace-Next $d "He's the king of the kitchen" "A friend so true and keen" $pos
$pos = Replace-Next $d "With a love for Nasi Goreng" "In his heart, love is seen" $pos
$pos = Replace-Next $d "His meals are always enrichin'" "Yaseen, Yaseen" $pos
$pos = Replace-Next $d "(Bridge)" "Bridge:" $pos
$pos = Replace-Next $d "With a smile on his face" "No matter the distance or time apart" $pos
$pos = Replace-Next $d "And a fork in his hand" "Yaseen will always have a place in my heart" $pos
$pos = Replace-Next $d "Sajeev's love for Nasi Goreng" "With his loyalty and friendship, I am blessed" $pos
$pos = Replace-Next $d "Is something truly grand" "Forever grateful for his love and zest" $pos
$pos = Replace-Next $d "(Chorus)" "Chorus:" $pos
$pos = Replace-Next $d "Sajeev, oh Sajeev" "Yaseen, Yaseen" $pos
$pos = Replace-Next $d "He's our Malaysian friend" "A friend so true and keen" $pos
$pos = Replace-Next $d "With a love for Nasi Goreng" "In his heart, love is seen" $pos
$pos = Replace-Next $d "That will never end" "Yaseen, Yaseen" $pos
$pos = Replace-Next $d "(Outro)" "Outro:" $pos
$pos = Replace-Next $d "So here's to Sajeev" "To my dear Pakistani friend" $pos
$pos = Replace-Next $d "Our food-loving mate" "Yaseen, may our bond never end" $pos
$pos = Replace-Next $d "With his love for Nasi Goreng" "In this world or the next, I'll always be" $pos
$pos = Replace-Next $d "He'll always dominate" "Grateful for your friendship eternally." $pos

Write-Host $d.Content.Text
